$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6540
$ws.Range("D27").Value = 6098568
$ws.Range("E27").Value = 932.5027522935779
$ws.Range("F27").Value = 9.915966386554631
$ws.Range("H27").Value = 25.08582184375903
